$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last data row (row 36) into a new row 37, mirroring the
# existing pattern of repeated Date / price entries already present in
# rows 33-36. Copy/PasteSpecial brings along both the values (so the
# same shared strings are reused) and the cell formatting/style.
$ws.Range("A36:B36").Copy()
$ws.Range("A37").PasteSpecial()
